# Auto-generated update script for Jenova Profits market-data refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 7371.643
$ws.Range("I40").Value = 5438
$ws.Range("K40").Value = 5438
$ws.Range("M40").Value = -5263
$ws.Range("H58").Value = 5309.591
$ws.Range("I58").Value = 2100.1538
$ws.Range("K58").Value = 6300.4614
$ws.Range("M58").Value = -6150.4614
$ws.Range("H62").Value = 8933251
$ws.Range("I62").Value = 15627941
$ws.Range("J62").Value = 6998
$ws.Range("K62").Value = 15627941
$ws.Range("L62").Value = 6998
$ws.Range("M62").Value = -15627317
$ws.Range("N62").Value = -8246
$ws.Range("H65").Value = 8933251
$ws.Range("I65").Value = 15627941
$ws.Range("J65").Value = 6998
$ws.Range("K65").Value = 78139705
$ws.Range("L65").Value = 34990
$ws.Range("M65").Value = -78136585
$ws.Range("N65").Value = -41230
$ws.Range("H70").Value = 335166.66
$ws.Range("I70").Value = 2000
$ws.Range("J70").Value = 501750
$ws.Range("K70").Value = 6000
$ws.Range("L70").Value = 1505250
$ws.Range("M70").Value = -5730
$ws.Range("N70").Value = -1505790
$ws.Range("H73").Value = 335166.66
$ws.Range("I73").Value = 2000
$ws.Range("J73").Value = 501750
$ws.Range("K73").Value = 6000
$ws.Range("L73").Value = 1505250
$ws.Range("M73").Value = -5064
$ws.Range("N73").Value = -1507122
$ws.Range("H106").Value = 2556.75
$ws.Range("I106").Value = 2394.2
$ws.Range("K106").Value = 2394.2
$ws.Range("M106").Value = -1763.2
$ws.Range("H132").Value = 2549.879
$ws.Range("I132").Value = 2349.6155
$ws.Range("K132").Value = 7048.8465
$ws.Range("M132").Value = -4518.8465
$ws.Range("H137").Value = 2391.0833
$ws.Range("I137").Value = 1945.6
$ws.Range("J137").Value = 4618.5
$ws.Range("K137").Value = 5836.799999999999
$ws.Range("L137").Value = 13855.5
$ws.Range("M137").Value = -3286.799999999999
$ws.Range("N137").Value = -18955.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 4956.6665
$ws.Range("I31").Value = 4956.6665
$ws.Range("K31").Value = 4956.6665
$ws.Range("M31").Value = -4662.6665
$ws.Range("H32").Value = 1467.37
$ws.Range("I32").Value = 1467.37
$ws.Range("K32").Value = 1467.37
$ws.Range("M32").Value = -1180.37
$ws.Range("H74").Value = 1622.1538
$ws.Range("I74").Value = 1627.4286
$ws.Range("J74").Value = 1600
$ws.Range("K74").Value = 1627.4286
$ws.Range("L74").Value = 1600
$ws.Range("M74").Value = -753.4286
$ws.Range("N74").Value = -3348
$ws.Range("H77").Value = 1622.1538
$ws.Range("I77").Value = 1627.4286
$ws.Range("J77").Value = 1600
$ws.Range("K77").Value = 8137.143
$ws.Range("L77").Value = 8000
$ws.Range("M77").Value = -3769.143
$ws.Range("N77").Value = -16736
$ws.Range("H97").Value = 712.7406999999999
$ws.Range("I97").Value = 712.7406999999999
$ws.Range("K97").Value = 712.7406999999999
$ws.Range("M97").Value = -216.7406999999999
$ws.Range("H102").Value = 1994.3846
$ws.Range("I102").Value = 1702.909
$ws.Range("J102").Value = 3597.5
$ws.Range("K102").Value = 1702.909
$ws.Range("L102").Value = 3597.5
$ws.Range("M102").Value = -80.90900000000011
$ws.Range("N102").Value = -6841.5
$ws.Range("H132").Value = 4680.136
$ws.Range("I132").Value = 4840.2104
$ws.Range("K132").Value = 14520.6312
$ws.Range("M132").Value = -11990.6312

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H40").Value = 54998.332
$ws.Range("J40").Value = 54998.332
$ws.Range("L40").Value = 54998.332
$ws.Range("N40").Value = -55528.332
$ws.Range("H134").Value = 46204.28
$ws.Range("I134").Value = 6005.095
$ws.Range("K134").Value = 18015.285
$ws.Range("M134").Value = -15480.285

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 1092.6666
$ws.Range("I94").Value = 637.3333
$ws.Range("K94").Value = 637.3333
$ws.Range("M94").Value = -186.3333
$ws.Range("H132").Value = 2578.1875
$ws.Range("I132").Value = 2352.1667
$ws.Range("K132").Value = 7056.500100000001
$ws.Range("M132").Value = -4526.500100000001
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 58747.5
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 58747.5
$ws.Range("M137").ClearContents()
$ws.Range("N137").Value = -68947.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3975
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()
$ws.Range("H12").Value = 261.5
$ws.Range("I12").Value = 125
$ws.Range("J12").Value = 281
$ws.Range("K12").Value = 375
$ws.Range("L12").Value = 843
$ws.Range("M12").Value = -202
$ws.Range("N12").Value = -1189
$ws.Range("H87").Value = 30014
$ws.Range("I87").Value = 30014
$ws.Range("K87").Value = 90042
$ws.Range("M87").Value = -88794
$ws.Range("H90").Value = 30014
$ws.Range("I90").Value = 30014
$ws.Range("K90").Value = 270126
$ws.Range("M90").Value = -263886
$ws.Range("H132").Value = 554985.6
$ws.Range("I132").Value = 144009.28
$ws.Range("J132").Value = 776280.6
$ws.Range("K132").Value = 1296083.52
$ws.Range("L132").Value = 6986525.399999999
$ws.Range("M132").Value = -1293553.52
$ws.Range("N132").Value = -6991585.399999999
$ws.Range("H133").Value = 8346.286
$ws.Range("I133").Value = 8969.799999999999
$ws.Range("K133").Value = 26909.4
$ws.Range("M133").Value = -21849.4
$ws.Range("H134").Value = 2001
$ws.Range("I134").Value = 2001
$ws.Range("K134").Value = 6003
$ws.Range("M134").Value = -933
$ws.Range("H136").Value = 951.6667
$ws.Range("I136").Value = 951.6667
$ws.Range("K136").Value = 2855.0001
$ws.Range("M136").Value = 2244.9999
$ws.Range("H137").Value = 2783.5557
$ws.Range("I137").Value = 2783.5557
$ws.Range("K137").Value = 8350.667099999999
$ws.Range("M137").Value = -3250.667099999999
$ws.Range("H138").Value = 4207.1816
$ws.Range("I138").Value = 4098.3335
$ws.Range("K138").Value = 12295.0005
$ws.Range("M138").Value = -7155.000499999998
$ws.Range("H139").Value = 7999.9287
$ws.Range("I139").Value = 7999.5
$ws.Range("K139").Value = 23998.5
$ws.Range("M139").Value = -18858.5
$ws.Range("H140").Value = 4341.2856
$ws.Range("I140").Value = 4066.6
$ws.Range("J140").Value = 5028
$ws.Range("K140").Value = 12199.8
$ws.Range("L140").Value = 15084
$ws.Range("M140").Value = -7019.799999999999
$ws.Range("N140").Value = -25444
$ws.Range("H141").Value = 4255.5713
$ws.Range("I141").Value = 4255.5713
$ws.Range("K141").Value = 12766.7139
$ws.Range("M141").Value = -7586.713899999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1956.1111
$ws.Range("I102").Value = 1281.7142
$ws.Range("K102").Value = 1281.7142
$ws.Range("M102").Value = 340.2858000000001
$ws.Range("H113").Value = 627860.8
$ws.Range("J113").Value = 4750
$ws.Range("L113").Value = 4750
$ws.Range("N113").Value = -9090
$ws.Range("H122").Value = 4007.5
$ws.Range("I122").Value = 1512.8572
$ws.Range("J122").Value = 7500
$ws.Range("K122").Value = 4538.571599999999
$ws.Range("L122").Value = 22500
$ws.Range("M122").Value = -2088.571599999999
$ws.Range("N122").Value = -27400
$ws.Range("H132").Value = 213956
$ws.Range("I132").Value = 17445
$ws.Range("J132").Value = 1000000
$ws.Range("K132").Value = 52335
$ws.Range("L132").Value = 3000000
$ws.Range("M132").Value = -49805
$ws.Range("N132").Value = -3005060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8567
$ws.Range("J7").Value = 6521.8335
$ws.Range("L7").Value = 6521.8335
$ws.Range("N7").Value = -6745.8335
$ws.Range("H40").Value = 3251.4814
$ws.Range("I40").Value = 3295.8076
$ws.Range("K40").Value = 3295.8076
$ws.Range("M40").Value = -3159.8076
$ws.Range("H45").Value = 1000000
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 1000000
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 1000000
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -1000814
$ws.Range("H55").Value = 1289.7142
$ws.Range("I55").Value = 143
$ws.Range("K55").Value = 143
$ws.Range("M55").Value = 30
$ws.Range("H82").Value = 4744.6665
$ws.Range("I82").Value = 4740
$ws.Range("K82").Value = 4740
$ws.Range("M82").Value = -4379
$ws.Range("H85").Value = 4744.6665
$ws.Range("I85").Value = 4740
$ws.Range("K85").Value = 4740
$ws.Range("M85").Value = -3492
$ws.Range("H126").Value = 8567
$ws.Range("J126").Value = 6521.8335
$ws.Range("L126").Value = 19565.5005
$ws.Range("N126").Value = -24505.5005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1129.25
$ws.Range("I107").Value = 1380.3043
$ws.Range("K107").Value = 4140.9129
$ws.Range("M107").Value = -2220.9129

Write-Host "Applied Jenova Profits market-data refresh."